$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.700.48"
$ws.Range("E2").Value = "  +3.51%  "
$ws.Range("D3").Value = "1.864.76"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.48%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +11.71%  "
$ws.Range("E9").Value = "  +7.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0697"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.38%  "
$ws.Range("E11").Value = "  +4.01%  "
$ws.Range("D12").Value = "2.134.46"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.14%  "
$ws.Range("D14").Value = "1.870.04"
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.682"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.32%  "
$ws.Range("D17").Value = "35.710.83"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").Value = "0.0₃0806"
$ws.Range("E20").Value = "  +4.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.96%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.43%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "3.347.23"
$ws.Range("E31").Value = "  +37.76%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0549"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.13%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.21%  "
$ws.Range("E34").Value = "  +4.45%  "
$ws.Range("E35").Value = "  +4.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "101.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +23.88%  "
$ws.Range("E37").Value = "  +7.52%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.370.56"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.05%  "
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("E41").Value = "  +5.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.89%  "
$ws.Range("E44").Value = "  +4.19%  "
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.88%  "
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D49").Value = "2.032.45"
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("E51").Value = "  +0.37%  "
